# Casos de prueba.xlsx - agregar casos de prueba #13 a #16
# (Crear, Ver, Editar y Borrar productos)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# --- Caso #13: Crear Producto -------------------------------------------
$ws.Range("A19").Value = "Caso #13"
$ws.Range("B19").Value = "Crear Producto"
$ws.Range("C19").Value = "Crear Producto "
$ws.Range("D19").Value = 45274
$ws.Range("E19").Value = "SI"
$ws.Range("F19").Value = "-"
$ws.Range("G19").Value = "OK"

# --- Caso #14: Ver Productos ---------------------------------------------
$ws.Range("A20").Value = "Caso #14"
$ws.Range("B20").Value = "Ver Productos"
$ws.Range("C20").Value = "Ver lista de productos"
$ws.Range("D20").Value = 45274
$ws.Range("E20").Value = "SI"
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = "OK"

# --- Caso #15 / Caso #16 id column filled first, matching original order -
$ws.Range("A21").Value = "Caso #15"
$ws.Range("A22").Value = "Caso #16"

# --- Caso #15: Editar producto --------------------------------------------
$ws.Range("B21").Value = "Editar porducto"
$ws.Range("C21").Value = "Editar producto (solo admin)"
$ws.Range("D21").Value = 45274
$ws.Range("E21").Value = "SI"
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = "OK"

# --- Caso #16: Borrar producto ---------------------------------------------
$ws.Range("B22").Value = "Borrar porducto"
$ws.Range("C22").Value = "Borrar producto (solo admin)"
$ws.Range("D22").Value = 45274
$ws.Range("E22").Value = "SI"
$ws.Range("F22").Value = "-"
$ws.Range("G22").Value = "OK"

# --- view state (best effort; scroll position) -----------------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F26").Select()
